# Examples from 11 to 19 added
# Adds a red-highlighted "INCONSISTENTE" label in A1 and an
# "ESEMPIO INCONSISTENTE" label in M8 (new row 8), widens column A and M
# to fit the new labels, and leaves the selection on M9 (just below the
# new note), matching the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "INCONSISTENTE" flag in A1, red fill ---
$ws.Range("A1").Value = "INCONSISTENTE"
$ws.Range("A1").Interior.Color = 255   # pure red (RGB 255,0,0)

# --- New "ESEMPIO INCONSISTENTE" note in M8, red fill ---
$ws.Range("M8").Value = "ESEMPIO INCONSISTENTE"
$ws.Range("M8").Interior.Color = 255   # pure red (RGB 255,0,0)

# --- Widen column A and M so the new labels are readable ---
$ws.Columns("A").ColumnWidth = 18.1
$ws.Columns("M").ColumnWidth = 22.6

# --- Leave the selection where the author left it ---
$ws.Range("M9").Select()
